# This script updates the three-digit x one-digit multiplication problems
# in the table to a new generated set of equations, per the commit
# "Update master to output generated at c8c62b6".
#
# Replacements are applied in document order using Find/Replace so that the
# one pair whose new text equals another pair's old text
# (547x7=3829 -> 717x2=1434, and 904x7=6328 -> 547x7=3829) is processed
# safely: the occurrence of "547x7=3829" is replaced away before the new
# "547x7=3829" is introduced later in the document.

$d = $word.ActiveDocument

$replacements = @(
    @("925×3=2775", "149×5=745"),
    @("101×8=808",  "202×5=1010"),
    @("705×8=5640", "430×6=2580"),
    @("205×9=1845", "919×8=7352"),
    @("687×5=3435", "339×8=2712"),
    @("610×4=2440", "328×7=2296"),
    @("279×9=2511", "527×8=4216"),
    @("160×2=320",  "917×9=8253"),
    @("301×9=2709", "208×2=416"),
    @("923×3=2769", "204×9=1836"),
    @("268×6=1608", "260×2=520"),
    @("871×6=5226", "929×7=6503"),
    @("149×6=894",  "212×2=424"),
    @("984×2=1968", "539×9=4851"),
    @("650×4=2600", "396×6=2376"),
    @("839×6=5034", "666×9=5994"),
    @("447×9=4023", "637×2=1274"),
    @("445×6=2670", "650×7=4550"),
    @("547×7=3829", "717×2=1434"),
    @("927×6=5562", "246×6=1476"),
    @("644×8=5152", "103×6=618"),
    @("904×7=6328", "547×7=3829"),
    @("705×4=2820", "712×7=4984"),
    @("690×7=4830", "272×4=1088"),
    @("541×3=1623", "902×5=4510")
)

foreach ($pair in $replacements) {
    $oldText = $pair[0]
    $newText = $pair[1]
    $found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                                      $true, 1, $false, $newText, 2)
    if (-not $found) {
        Write-Output "WARNING: could not find '$oldText'"
    }
}

Write-Output "Done."
